{"js": "// 1) Insert a new \"Meta description\" paragraph right after the first\n//    (Heading1) paragraph: \"Play Age of the Gods Goddess of Wisdom for Free | Review\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst heading = paragraphs.items[0];\n\n// Create a placeholder empty paragraph right after the heading, then\n// overwrite it with precise OOXML so the run layout (empty lead run,\n// bold \"Meta description\" run, plain trailing run) matches exactly.\nconst metaPara = heading.insertParagraph(\"\", \"After\");\nawait context.sync();\n\nconst metaOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Find out everything you need to know about Age of the Gods Goddess of Wisdom, a Playtech online slot game with a Greek mythology theme. Play it for free and hit any of the 4 progressive jackpots at any time.</w:t></w:r></w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\nmetaPara.insertOoxml(metaOoxml, \"Replace\");\nawait context.sync();\n\n// 2) Near the end of the document: remove the duplicated bold\n//    \"Play Age of the Gods Goddess of Wisdom for Free | Review\" paragraph,\n//    and rewrite the final (italic) paragraph's text to the new image prompt.\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst count = items.length;\nconst lastPara = items[count - 1];\nconst secondLastPara = items[count - 2];\n\nif (\n  secondLastPara.text === \"Play Age of the Gods Goddess of Wisdom for Free | Review\" &&\n  lastPara.text ===\n    \"Find out everything you need to know about Age of the Gods Goddess of Wisdom, a Playtech online slot game with a Greek mythology theme. Play it for free and hit any of the 4 progressive jackpots at any time.\"\n) {\n  secondLastPara.delete();\n  await context.sync();\n\n  lastPara.insertText(\n    'Create a feature image for \"Age of the Gods: Goddess of Wisdom\". The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be standing in front of an ancient Greek temple while holding a shield with the game title written on it. Athena, the central figure of the game, should be standing next to the warrior with a confident stance. The symbols of the game, including the Gorgoneion, Olive Branches, Helmets, and Armor, should be seen floating around the two figures. The image should be colorful and eye-catching to attract players\\' attention.',\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# 1) Insert a new \"Meta description\" paragraph right after the first\n#    (Heading1) paragraph: \"Play Age of the Gods Goddess of Wisdom for Free | Review\".\n$d = $word.ActiveDocument\n\n$heading = $d.Paragraphs(1)\n$headingRange = $heading.Range\n$headingRange.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs(2)\n$metaRange = $metaPara.Range\n\n$metaOoxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Find out everything you need to know about Age of the Gods Goddess of Wisdom, a Playtech online slot game with a Greek mythology theme. Play it for free and hit any of the 4 progressive jackpots at any time.</w:t></w:r></w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$metaRange.InsertXML($metaOoxml)\n\n# 2) Near the end of the document: remove the duplicated bold\n#    \"Play Age of the Gods Goddess of Wisdom for Free | Review\" paragraph,\n#    and rewrite the final (italic) paragraph's text to the new image prompt.\n$count = $d.Paragraphs.Count\n$secondLastPara = $d.Paragraphs($count - 1)\n$lastPara = $d.Paragraphs($count)\n\nif (($secondLastPara.Range.Text -eq \"Play Age of the Gods Goddess of Wisdom for Free | Review`r\") -and `\n    ($lastPara.Range.Text -eq \"Find out everything you need to know about Age of the Gods Goddess of Wisdom, a Playtech online slot game with a Greek mythology theme. Play it for free and hit any of the 4 progressive jackpots at any time.`r\")) {\n\n    $secondLastPara.Range.Delete()\n\n    # Rewrite just the text content of the final paragraph (excluding its\n    # trailing paragraph mark) so the existing run formatting (italic) is\n    # preserved and no smart-quote autocorrect is triggered (unlike\n    # Find/Replace, a direct Range.Text assignment leaves punctuation as-is).\n    $newCount = $d.Paragraphs.Count\n    $finalPara = $d.Paragraphs($newCount)\n    $finalRange = $finalPara.Range\n    $textOnlyRange = $d.Range($finalRange.Start, $finalRange.End - 1)\n    $textOnlyRange.Text = \"Create a feature image for \"\"Age of the Gods: Goddess of Wisdom\"\". The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be standing in front of an ancient Greek temple while holding a shield with the game title written on it. Athena, the central figure of the game, should be standing next to the warrior with a confident stance. The symbols of the game, including the Gorgoneion, Olive Branches, Helmets, and Armor, should be seen floating around the two figures. The image should be colorful and eye-catching to attract players' attention.\"\n}\n"}
